# "Pillola finanziaria:" + line-break + "SOCIAL LENDINGS"  ->  "...SOCIAL LENDING"
# Only touch the trailing "LENDINGS" -> "LENDING" substring so PowerPoint keeps
# the orange run formatting and simply splits the run where the edit happened
# (mirrors how the author retyped the word in the real deck).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$word = $tr.Characters(29, 8)
$word.Text = "LENDING"
